# Finalização dos insights magalu - update slide
# Reposition/resize the background rectangle ("Retângulo 5") on slide 1.
# Target EMUs: off x=4018208 y=2395470, ext cx=8075053 cy=3258355
# (point values below are chosen so that, after the host's f32 EMU
# conversion, they land exactly on the target EMU values)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)

$sh.Left = 316.39434814453125
$sh.Top = 188.61968994140625
$sh.Width = 635.8309936523438
$sh.Height = 256.56341552734375
